$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.665.47'
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").Value = '3.285.12'
$ws.Range("E3").Value = '  +5.10%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.03'
$ws.Range("E5").Value = '  +2.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.23'
$ws.Range("E6").Value = '  +3.85%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.282.79'
$ws.Range("E8").Value = '  +5.11%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  +3.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  +4.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +3.33%  '
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.40'
$ws.Range("E14").Value = '  +1.50%  '
$ws.Range("D15").Value = '3.823.56'
$ws.Range("E15").Value = '  +5.06%  '
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '3.283.00'
$ws.Range("D18").Value = '63.716.87'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.90'
$ws.Range("E20").Value = '  +1.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  +4.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.99'
$ws.Range("E23").Value = '  +4.66%  '
$ws.Range("E24").Value = '  +4.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.11'
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +2.27%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("E29").Value = '  +6.40%  '
$ws.Range("E30").Value = '  +2.94%  '
$ws.Range("E31").Value = '  +3.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.61'
$ws.Range("E32").Value = '  +7.94%  '
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +3.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("E36").Value = '  +3.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.27'
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("D38").Value = '0.0₃0734'
$ws.Range("E38").Value = '  +8.54%  '
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '427.23'
$ws.Range("E40").Value = '  +3.35%  '
$ws.Range("D41").Value = '3.052.33'
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.31'
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.72'
$ws.Range("E43").Value = '  +1.90%  '
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E46").Value = '  +3.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.10'
$ws.Range("E48").Value = '  +3.52%  '
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.86'
$ws.Range("E50").Value = '  +3.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.27'
$ws.Range("E51").Value = '  +1.30%  '
